# Updates TPM-derived NATMI ligand-receptor metrics for Ltf-Lrp11 (YoungD4) with new TPM values.
# Ligand-expressing cell count for Ltf dropped from 2 -> 1 in each sending cluster, which
# cascades through detection rate, average/total expression, derived specificities and edge metrics.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "E2" = 1
    "F2" = 0.3333333333333333
    "G2" = 0.02542733333333333
    "H2" = 0.076282
    "I2" = 0.3241752404264994
    "J2" = 0.3241752404264994
    "M2" = 0.9772823333333333
    "N2" = 2.931847
    "O2" = 0.1949382298804705
    "P2" = 0.1949382298804704
    "Q2" = 0.02484968365044444
    "R2" = 0.223647152854
    "S2" = 0.06319414753981771
    "T2" = 0.06319414753981771
    "E3" = 1
    "F3" = 0.3333333333333333
    "G3" = 0.02542733333333333
    "H3" = 0.076282
    "I3" = 0.3241752404264994
    "J3" = 0.3241752404264994
    "O3" = 0.304339058092031
    "P3" = 0.3043390580920309
    "Q3" = 0.03879551651155556
    "R3" = 0.349159648604
    "S3" = 0.09865918732815851
    "T3" = 0.0986591873281585
    "E4" = 1
    "F4" = 0.3333333333333333
    "G4" = 0.02542733333333333
    "H4" = 0.076282
    "I4" = 0.3241752404264994
    "J4" = 0.3241752404264994
    "M4" = 0.741802
    "N4" = 2.225406
    "O4" = 0.147967034570828
    "P4" = 0.1479670345708279
    "Q4" = 0.01886204672133333
    "R4" = 0.169758420492
    "S4" = 0.0479672490071943
    "T4" = 0.0479672490071943
    "E5" = 1
    "F5" = 0.3333333333333333
    "G5" = 0.02542733333333333
    "H5" = 0.076282
    "I5" = 0.3241752404264994
    "J5" = 0.3241752404264994
    "M5" = 1.021925666666667
    "N5" = 3.065777
    "O5" = 0.2038432229199747
    "P5" = 0.2038432229199746
    "Q5" = 0.02598484456822222
    "R5" = 0.233863601114
    "S5" = 0.0660809257993953
    "T5" = 0.0660809257993953
    "E6" = 1
    "F6" = 0.3333333333333333
    "G6" = 0.02542733333333333
    "H6" = 0.076282
    "I6" = 0.3241752404264994
    "J6" = 0.3241752404264994
    "K6" = 3
    "L6" = 1
    "M6" = 0.7465416666666668
    "N6" = 2.239625
    "O6" = 0.148912454536696
    "P6" = 0.148912454536696
    "Q6" = 0.01898256380555556
    "R6" = 0.17084307425
    "S6" = 0.0482737307519336
    "T6" = 0.0482737307519336
    "E7" = 1
    "F7" = 0.3333333333333333
    "G7" = 0.05300966666666667
    "H7" = 0.159029
    "I7" = 0.6758247595735006
    "J7" = 0.6758247595735006
    "M7" = 0.9772823333333333
    "N7" = 2.931847
    "O7" = 0.1949382298804705
    "P7" = 0.1949382298804704
    "Q7" = 0.05180541072922222
    "R7" = 0.466248696563
    "S7" = 0.1317440823406527
    "T7" = 0.1317440823406527
    "E8" = 1
    "F8" = 0.3333333333333333
    "G8" = 0.05300966666666667
    "H8" = 0.159029
    "I8" = 0.6758247595735006
    "J8" = 0.6758247595735006
    "O8" = 0.304339058092031
    "P8" = 0.3043390580920309
    "Q8" = 0.08087900415977778
    "R8" = 0.727911037438
    "S8" = 0.2056798707638725
    "T8" = 0.2056798707638725
    "E9" = 1
    "F9" = 0.3333333333333333
    "G9" = 0.05300966666666667
    "H9" = 0.159029
    "I9" = 0.6758247595735006
    "J9" = 0.6758247595735006
    "M9" = 0.741802
    "N9" = 2.225406
    "O9" = 0.147967034570828
    "P9" = 0.1479670345708279
    "Q9" = 0.03932267675266667
    "R9" = 0.353904090774
    "S9" = 0.09999978556363366
    "T9" = 0.09999978556363365
    "E10" = 1
    "F10" = 0.3333333333333333
    "G10" = 0.05300966666666667
    "H10" = 0.159029
    "I10" = 0.6758247595735006
    "J10" = 0.6758247595735006
    "M10" = 1.021925666666667
    "N10" = 3.065777
    "O10" = 0.2038432229199747
    "P10" = 0.2038432229199746
    "Q10" = 0.05417193894811111
    "R10" = 0.487547450533
    "S10" = 0.1377622971205794
    "T10" = 0.1377622971205794
    "E11" = 1
    "F11" = 0.3333333333333333
    "G11" = 0.05300966666666667
    "H11" = 0.159029
    "I11" = 0.6758247595735006
    "J11" = 0.6758247595735006
    "K11" = 3
    "L11" = 1
    "M11" = 0.7465416666666668
    "N11" = 2.239625
    "O11" = 0.148912454536696
    "P11" = 0.148912454536696
    "Q11" = 0.03957392490277779
    "R11" = 0.356165324125
    "S11" = 0.1006387237847624
    "T11" = 0.1006387237847624
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
